$d = $word.ActiveDocument

# 1. "Incluir títul" + _GoBack bookmark + "os a receber" -> merge back into a
#    single run "Incluir títulos a receber" (removes the old bookmark split).
$d.Content.Find.Execute("Incluir títulos a receber", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Incluir títulos a receber", 2)

# 2. "Gerenciar Tipos de Despesa" -> "Gerenciar Tipos de Pagamento"
$d.Content.Find.Execute("Gerenciar Tipos de Despesa", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Gerenciar Tipos de Pagamento", 2)

# 3. "Quitar contas a receber" -> split into "Quitar contas a re" + _GoBack
#    bookmark + "ceber" (simulating the cursor being left there after typing).
$r = $d.Content
$r.Find.Execute("Quitar contas a re")
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r)
